# Update Cxcl12-Itgb1 LR-pair data following Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 109.026058
$ws.Cells.Item(2, 8).Value = 327.078174
$ws.Cells.Item(2, 9).Value = 0.3049840938689738
$ws.Cells.Item(2, 10).Value = 0.3049840938689738
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 117.044563
$ws.Cells.Item(2, 14).Value = 351.133689
$ws.Cells.Item(2, 15).Value = 0.3245365645427815
$ws.Cells.Item(2, 16).Value = 0.3245365645427815
$ws.Cells.Item(2, 17).Value = 12760.90731422265
$ws.Cells.Item(2, 18).Value = 114848.1658280039
$ws.Cells.Item(2, 19).Value = 0.09897849006442994
$ws.Cells.Item(2, 20).Value = 0.09897849006442994

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 109.026058
$ws.Cells.Item(3, 8).Value = 327.078174
$ws.Cells.Item(3, 9).Value = 0.3049840938689738
$ws.Cells.Item(3, 10).Value = 0.3049840938689738
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 101.5800373333333
$ws.Cells.Item(3, 14).Value = 304.740112
$ws.Cells.Item(3, 15).Value = 0.281657135515876
$ws.Cells.Item(3, 16).Value = 0.281657135515876
$ws.Cells.Item(3, 17).Value = 11074.87104194616
$ws.Cells.Item(3, 18).Value = 99673.83937751548
$ws.Cells.Item(3, 19).Value = 0.0859009462570402
$ws.Cells.Item(3, 20).Value = 0.08590094625704019

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 109.026058
$ws.Cells.Item(4, 8).Value = 327.078174
$ws.Cells.Item(4, 9).Value = 0.3049840938689738
$ws.Cells.Item(4, 10).Value = 0.3049840938689738
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 142.0267893333333
$ws.Cells.Item(4, 14).Value = 426.080368
$ws.Cells.Item(4, 15).Value = 0.3938062999413425
$ws.Cells.Item(4, 16).Value = 0.3938062999413425
$ws.Cells.Item(4, 17).Value = 15484.62097140978
$ws.Cells.Item(4, 18).Value = 139361.588742688
$ws.Cells.Item(4, 19).Value = 0.1201046575475037
$ws.Cells.Item(4, 20).Value = 0.1201046575475037

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 89.97721833333333
$ws.Cells.Item(5, 8).Value = 269.931655
$ws.Cells.Item(5, 9).Value = 0.2516978134001918
$ws.Cells.Item(5, 10).Value = 0.2516978134001917
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 117.044563
$ws.Cells.Item(5, 14).Value = 351.133689
$ws.Cells.Item(5, 15).Value = 0.3245365645427815
$ws.Cells.Item(5, 16).Value = 0.3245365645427815
$ws.Cells.Item(5, 17).Value = 10531.34419978059
$ws.Cells.Item(5, 18).Value = 94782.0977980253
$ws.Cells.Item(5, 19).Value = 0.0816851436638283
$ws.Cells.Item(5, 20).Value = 0.08168514366382829

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 89.97721833333333
$ws.Cells.Item(6, 8).Value = 269.931655
$ws.Cells.Item(6, 9).Value = 0.2516978134001918
$ws.Cells.Item(6, 10).Value = 0.2516978134001917
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 101.5800373333333
$ws.Cells.Item(6, 14).Value = 304.740112
$ws.Cells.Item(6, 15).Value = 0.281657135515876
$ws.Cells.Item(6, 16).Value = 0.281657135515876
$ws.Cells.Item(6, 17).Value = 9139.889197449484
$ws.Cells.Item(6, 18).Value = 82259.00277704536
$ws.Cells.Item(6, 19).Value = 0.07089248513790748
$ws.Cells.Item(6, 20).Value = 0.07089248513790745

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 89.97721833333333
$ws.Cells.Item(7, 8).Value = 269.931655
$ws.Cells.Item(7, 9).Value = 0.2516978134001918
$ws.Cells.Item(7, 10).Value = 0.2516978134001917
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 142.0267893333333
$ws.Cells.Item(7, 14).Value = 426.080368
$ws.Cells.Item(7, 15).Value = 0.3938062999413425
$ws.Cells.Item(7, 16).Value = 0.3938062999413425
$ws.Cells.Item(7, 17).Value = 12779.17543302767
$ws.Cells.Item(7, 18).Value = 115012.578897249
$ws.Cells.Item(7, 19).Value = 0.09912018459845598
$ws.Cells.Item(7, 20).Value = 0.09912018459845595

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 158.477852
$ws.Cells.Item(8, 8).Value = 475.433556
$ws.Cells.Item(8, 9).Value = 0.4433180927308344
$ws.Cells.Item(8, 10).Value = 0.4433180927308344
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 117.044563
$ws.Cells.Item(8, 14).Value = 351.133689
$ws.Cells.Item(8, 15).Value = 0.3245365645427815
$ws.Cells.Item(8, 16).Value = 0.3245365645427815
$ws.Cells.Item(8, 17).Value = 18548.97093251867
$ws.Cells.Item(8, 18).Value = 166940.7383926681
$ws.Cells.Item(8, 19).Value = 0.1438729308145232
$ws.Cells.Item(8, 20).Value = 0.1438729308145232

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 158.477852
$ws.Cells.Item(9, 8).Value = 475.433556
$ws.Cells.Item(9, 9).Value = 0.4433180927308344
$ws.Cells.Item(9, 10).Value = 0.4433180927308344
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 101.5800373333333
$ws.Cells.Item(9, 14).Value = 304.740112
$ws.Cells.Item(9, 15).Value = 0.281657135515876
$ws.Cells.Item(9, 16).Value = 0.281657135515876
$ws.Cells.Item(9, 17).Value = 16098.18612266647
$ws.Cells.Item(9, 18).Value = 144883.6751039982
$ws.Cells.Item(9, 19).Value = 0.1248637041209283
$ws.Cells.Item(9, 20).Value = 0.1248637041209283

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 158.477852
$ws.Cells.Item(10, 8).Value = 475.433556
$ws.Cells.Item(10, 9).Value = 0.4433180927308344
$ws.Cells.Item(10, 10).Value = 0.4433180927308344
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 142.0267893333333
$ws.Cells.Item(10, 14).Value = 426.080368
$ws.Cells.Item(10, 15).Value = 0.3938062999413425
$ws.Cells.Item(10, 16).Value = 0.3938062999413425
$ws.Cells.Item(10, 17).Value = 22508.10050000318
$ws.Cells.Item(10, 18).Value = 202572.9045000286
$ws.Cells.Item(10, 19).Value = 0.1745814577953829
$ws.Cells.Item(10, 20).Value = 0.1745814577953829

